$wb = $excel.ActiveWorkbook

$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")
$wsGlobal.Range("C2").Value = 5.219841746182212
$wsGlobal.Range("C3").Value = 12.202597394334816
$wsGlobal.Range("C4").Value = 0.6515745358190308
$wsGlobal.Range("C6").Value = 5.030926638536917
$wsGlobal.Range("C7").Value = 11.767636273709758
$wsGlobal.Range("C8").Value = 0.7154790360916665
$wsGlobal.Range("C10").Value = 5.030926638536917
$wsGlobal.Range("C11").Value = 11.767636273709758
$wsGlobal.Range("C12").Value = 0.7154790360916665
$wsGlobal.Range("C14").Value = 5.073941390864349
$wsGlobal.Range("C15").Value = 11.866674117262637
$wsGlobal.Range("C16").Value = 0.4632885255976896
$wsGlobal.Range("C18").Value = 5.064594300942521
$wsGlobal.Range("C19").Value = 11.845153230142188
$wsGlobal.Range("C20").Value = 0.6833463972097202

$wsLanding = $wb.Worksheets.Item("LANDING GEARS")
$wsLanding.Range("C2").Value = 12.298109362990228
